# Scheduled runner update: refresh market-price-derived profit columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ,
# i.e. columns H-N) for specific leve rows across several job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 5196.8
$ws.Cells.Item(33, 9).Value = 7822
$ws.Cells.Item(33, 10).Value = 321.42856
$ws.Cells.Item(33, 11).Value = 7822
$ws.Cells.Item(33, 12).Value = 321.42856
$ws.Cells.Item(33, 13).Value = -7593
$ws.Cells.Item(33, 14).Value = -779.4285600000001

$ws.Cells.Item(37, 8).Value = 868.5714
$ws.Cells.Item(37, 10).Value = 868.5714
$ws.Cells.Item(37, 12).Value = 2605.7142
$ws.Cells.Item(37, 14).Value = -2857.7142

$ws.Cells.Item(112, 8).Value = 1721.6666
$ws.Cells.Item(112, 9).Value = 1500
$ws.Cells.Item(112, 10).Value = 1741.8182
$ws.Cells.Item(112, 11).Value = 4500
$ws.Cells.Item(112, 12).Value = 5225.4546
$ws.Cells.Item(112, 13).Value = -3392
$ws.Cells.Item(112, 14).Value = -7441.4546

$ws.Cells.Item(113, 8).Value = 7420.8
$ws.Cells.Item(113, 9).Value = 2527.6667
$ws.Cells.Item(113, 10).Value = 11937.538
$ws.Cells.Item(113, 11).Value = 2527.6667
$ws.Cells.Item(113, 12).Value = 11937.538
$ws.Cells.Item(113, 13).Value = 726.3332999999998
$ws.Cells.Item(113, 14).Value = -18445.538

$ws.Cells.Item(116, 8).Value = 2303.5625
$ws.Cells.Item(116, 9).Value = 2104.3845
$ws.Cells.Item(116, 11).Value = 2104.3845
$ws.Cells.Item(116, 13).Value = 1337.6155

$ws.Cells.Item(129, 8).Value = 877.18866
$ws.Cells.Item(129, 10).Value = 895.88
$ws.Cells.Item(129, 12).Value = 2687.64
$ws.Cells.Item(129, 14).Value = -12687.64

$ws.Cells.Item(138, 8).Value = 3216.9348
$ws.Cells.Item(138, 9).Value = 1729.6875
$ws.Cells.Item(138, 10).Value = 4010.1333
$ws.Cells.Item(138, 11).Value = 5189.0625
$ws.Cells.Item(138, 12).Value = 12030.3999
$ws.Cells.Item(138, 13).Value = -49.0625
$ws.Cells.Item(138, 14).Value = -22310.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 13230.454
$ws.Cells.Item(63, 9).Value = 14253.5
$ws.Cells.Item(63, 10).Value = 3000
$ws.Cells.Item(63, 11).Value = 14253.5
$ws.Cells.Item(63, 12).Value = 3000
$ws.Cells.Item(63, 13).Value = -13567.5
$ws.Cells.Item(63, 14).Value = -4372

$ws.Cells.Item(66, 8).Value = 13230.454
$ws.Cells.Item(66, 9).Value = 14253.5
$ws.Cells.Item(66, 10).Value = 3000
$ws.Cells.Item(66, 11).Value = 71267.5
$ws.Cells.Item(66, 12).Value = 15000
$ws.Cells.Item(66, 13).Value = -67835.5
$ws.Cells.Item(66, 14).Value = -21864

$ws.Cells.Item(97, 8).Value = 1799.45
$ws.Cells.Item(97, 9).Value = 1554.3334
$ws.Cells.Item(97, 10).Value = 4005.5
$ws.Cells.Item(97, 11).Value = 1554.3334
$ws.Cells.Item(97, 12).Value = 4005.5
$ws.Cells.Item(97, 13).Value = -1058.3334
$ws.Cells.Item(97, 14).Value = -4997.5

$ws.Cells.Item(102, 8).Value = 2086.5557
$ws.Cells.Item(102, 9).Value = 2012.25
$ws.Cells.Item(102, 10).Value = 2235.1667
$ws.Cells.Item(102, 11).Value = 2012.25
$ws.Cells.Item(102, 12).Value = 2235.1667
$ws.Cells.Item(102, 13).Value = -390.25
$ws.Cells.Item(102, 14).Value = -5479.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1725.25
$ws.Cells.Item(86, 9).Value = 2107.4
$ws.Cells.Item(86, 10).Value = 1284.3077
$ws.Cells.Item(86, 11).Value = 2107.4
$ws.Cells.Item(86, 12).Value = 1284.3077
$ws.Cells.Item(86, 13).Value = -984.4000000000001
$ws.Cells.Item(86, 14).Value = -3530.3077

$ws.Cells.Item(89, 8).Value = 1725.25
$ws.Cells.Item(89, 9).Value = 2107.4
$ws.Cells.Item(89, 10).Value = 1284.3077
$ws.Cells.Item(89, 11).Value = 10537
$ws.Cells.Item(89, 12).Value = 6421.538500000001
$ws.Cells.Item(89, 13).Value = -4921
$ws.Cells.Item(89, 14).Value = -17653.5385

$ws.Cells.Item(94, 8).Value = 737.03845
$ws.Cells.Item(94, 9).Value = 702.86365
$ws.Cells.Item(94, 10).Value = 925
$ws.Cells.Item(94, 11).Value = 702.86365
$ws.Cells.Item(94, 12).Value = 925
$ws.Cells.Item(94, 13).Value = -251.86365
$ws.Cells.Item(94, 14).Value = -1827

$ws.Cells.Item(99, 8).Value = 1204.0938
$ws.Cells.Item(99, 9).Value = 883
$ws.Cells.Item(99, 10).Value = 2024.6666
$ws.Cells.Item(99, 11).Value = 883
$ws.Cells.Item(99, 12).Value = 2024.6666
$ws.Cells.Item(99, 13).Value = 615
$ws.Cells.Item(99, 14).Value = -5020.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 2045.875
$ws.Cells.Item(86, 9).Value = 1964.125
$ws.Cells.Item(86, 10).Value = 2291.125
$ws.Cells.Item(86, 11).Value = 1964.125
$ws.Cells.Item(86, 12).Value = 2291.125
$ws.Cells.Item(86, 13).Value = -841.125
$ws.Cells.Item(86, 14).Value = -4537.125

$ws.Cells.Item(89, 8).Value = 2045.875
$ws.Cells.Item(89, 9).Value = 1964.125
$ws.Cells.Item(89, 10).Value = 2291.125
$ws.Cells.Item(89, 11).Value = 9820.625
$ws.Cells.Item(89, 12).Value = 11455.625
$ws.Cells.Item(89, 13).Value = -4204.625
$ws.Cells.Item(89, 14).Value = -22687.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1297.48
$ws.Cells.Item(5, 9).Value = 687.35297
$ws.Cells.Item(5, 10).Value = 2594
$ws.Cells.Item(5, 11).Value = 2062.05891
$ws.Cells.Item(5, 12).Value = 7782
$ws.Cells.Item(5, 13).Value = -1950.05891
$ws.Cells.Item(5, 14).Value = -8006

$ws.Cells.Item(107, 8).Value = 683
$ws.Cells.Item(107, 9).Value = 683
$ws.Cells.Item(107, 11).Value = 2049
$ws.Cells.Item(107, 13).Value = -129

$ws.Cells.Item(119, 8).Value = 6517.5557
$ws.Cells.Item(119, 9).Value = 4900
$ws.Cells.Item(119, 11).Value = 14700
$ws.Cells.Item(119, 13).Value = -9862

$ws.Cells.Item(122, 8).Value = 741.28
$ws.Cells.Item(122, 9).Value = 321.27274
$ws.Cells.Item(122, 10).Value = 1071.2858
$ws.Cells.Item(122, 11).Value = 2891.45466
$ws.Cells.Item(122, 12).Value = 9641.572200000001
$ws.Cells.Item(122, 13).Value = -441.4546599999999
$ws.Cells.Item(122, 14).Value = -14541.5722

$ws.Cells.Item(132, 8).Value = 1001.13635
$ws.Cells.Item(132, 9).Value = 799.7857
$ws.Cells.Item(132, 10).Value = 1353.5
$ws.Cells.Item(132, 11).Value = 7198.071300000001
$ws.Cells.Item(132, 12).Value = 12181.5
$ws.Cells.Item(132, 13).Value = -4668.071300000001
$ws.Cells.Item(132, 14).Value = -17241.5

$ws.Cells.Item(135, 8).Value = 1297.48
$ws.Cells.Item(135, 9).Value = 687.35297
$ws.Cells.Item(135, 10).Value = 2594
$ws.Cells.Item(135, 11).Value = 6186.17673
$ws.Cells.Item(135, 12).Value = 23346
$ws.Cells.Item(135, 13).Value = -3651.17673
$ws.Cells.Item(135, 14).Value = -28416

$ws.Cells.Item(140, 8).Value = 1463.2963
$ws.Cells.Item(140, 9).Value = 1094.3125
$ws.Cells.Item(140, 11).Value = 3282.9375
$ws.Cells.Item(140, 13).Value = 1897.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5399
$ws.Cells.Item(62, 9).Value = 4997.5
$ws.Cells.Item(62, 10).Value = 5666.6665
$ws.Cells.Item(62, 11).Value = 4997.5
$ws.Cells.Item(62, 12).Value = 5666.6665
$ws.Cells.Item(62, 13).Value = -4373.5
$ws.Cells.Item(62, 14).Value = -6914.6665

$ws.Cells.Item(65, 8).Value = 5399
$ws.Cells.Item(65, 9).Value = 4997.5
$ws.Cells.Item(65, 10).Value = 5666.6665
$ws.Cells.Item(65, 11).Value = 24987.5
$ws.Cells.Item(65, 12).Value = 28333.3325
$ws.Cells.Item(65, 13).Value = -21867.5
$ws.Cells.Item(65, 14).Value = -34573.3325

$ws.Cells.Item(96, 8).Value = 1384
$ws.Cells.Item(96, 9).Value = 1246.2222
$ws.Cells.Item(96, 10).Value = 2004
$ws.Cells.Item(96, 11).Value = 1246.2222
$ws.Cells.Item(96, 12).Value = 2004
$ws.Cells.Item(96, 13).Value = 126.7778000000001
$ws.Cells.Item(96, 14).Value = -4750

$ws.Cells.Item(132, 8).Value = 1940.3715
$ws.Cells.Item(132, 9).Value = 1555.9259
$ws.Cells.Item(132, 11).Value = 4667.7777
$ws.Cells.Item(132, 13).Value = -2137.7777
